# Update the division-problem worksheet: replace the text in each
# populated table cell with the new "dividend÷divisor=" expression.
# The 20x5 table only has content in every 4th row (1,5,9,13,17); we
# address cells directly by (row,col) and overwrite Cell.Range.Text so
# existing run formatting (TimeNewRoman / sz 30) is preserved and there
# is no ambiguity from duplicate expressions (e.g. "67÷8=" appears
# twice in the original document with two different replacements).

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text  = "59÷7="
$t.Cell(1,2).Range.Text  = "57÷3="
$t.Cell(1,3).Range.Text  = "96÷6="
$t.Cell(1,4).Range.Text  = "17÷6="
$t.Cell(1,5).Range.Text  = "60÷9="

$t.Cell(5,1).Range.Text  = "68÷6="
$t.Cell(5,2).Range.Text  = "60÷3="
$t.Cell(5,3).Range.Text  = "72÷8="
$t.Cell(5,4).Range.Text  = "63÷6="
$t.Cell(5,5).Range.Text  = "95÷2="

$t.Cell(9,1).Range.Text  = "60÷5="
$t.Cell(9,2).Range.Text  = "21÷6="
$t.Cell(9,3).Range.Text  = "63÷5="
$t.Cell(9,4).Range.Text  = "22÷3="
$t.Cell(9,5).Range.Text  = "11÷5="

$t.Cell(13,1).Range.Text = "64÷5="
$t.Cell(13,2).Range.Text = "36÷2="
$t.Cell(13,3).Range.Text = "16÷9="
$t.Cell(13,4).Range.Text = "93÷3="
$t.Cell(13,5).Range.Text = "72÷9="

$t.Cell(17,1).Range.Text = "96÷4="
$t.Cell(17,2).Range.Text = "10÷4="
$t.Cell(17,3).Range.Text = "10÷2="
$t.Cell(17,4).Range.Text = "13÷2="
$t.Cell(17,5).Range.Text = "17÷3="

Write-Host "Updated 25 division cells."
